$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 751.80725  # H17
$ws.Cells.Item(17, 10).Value = 660.2941  # J17
$ws.Cells.Item(17, 12).Value = 1980.8823  # L17
$ws.Cells.Item(17, 14).Value = -2316.8823  # N17
$ws.Cells.Item(43, 8).Value = 910.625  # H43
$ws.Cells.Item(43, 9).Value = 756.4  # I43
$ws.Cells.Item(43, 10).Value = 980.7273  # J43
$ws.Cells.Item(43, 11).Value = 756.4  # K43
$ws.Cells.Item(43, 12).Value = 980.7273  # L43
$ws.Cells.Item(43, 13).Value = -687.4  # M43
$ws.Cells.Item(43, 14).Value = -1118.7273  # N43
$ws.Cells.Item(98, 8).Value = 2772.0435  # H98
$ws.Cells.Item(98, 9).Value = 1303  # I98
$ws.Cells.Item(98, 11).Value = 1303  # K98
$ws.Cells.Item(98, 13).Value = 195  # M98
$ws.Cells.Item(116, 8).Value = 485224.34  # H116
$ws.Cells.Item(116, 9).Value = 1252613.1  # I116
$ws.Cells.Item(116, 11).Value = 1252613.1  # K116
$ws.Cells.Item(116, 13).Value = -1249171.1  # M116
$ws.Cells.Item(122, 8).Value = 2772.0435  # H122
$ws.Cells.Item(122, 9).Value = 1303  # I122
$ws.Cells.Item(122, 11).Value = 3909  # K122
$ws.Cells.Item(122, 13).Value = -1459  # M122
$ws.Cells.Item(129, 8).Value = 833.8  # H129
$ws.Cells.Item(129, 10).Value = 969.7447  # J129
$ws.Cells.Item(129, 12).Value = 2909.2341  # L129
$ws.Cells.Item(129, 14).Value = -12909.2341  # N129
$ws.Cells.Item(135, 8).Value = 398.80646  # H135
$ws.Cells.Item(135, 9).Value = 368.6  # I135
$ws.Cells.Item(135, 11).Value = 3317.4  # K135
$ws.Cells.Item(135, 13).Value = -782.4000000000001  # M135
$ws.Cells.Item(137, 8).Value = 3887.9285  # H137
$ws.Cells.Item(137, 10).Value = 4321.8096  # J137
$ws.Cells.Item(137, 12).Value = 12965.4288  # L137
$ws.Cells.Item(137, 14).Value = -18065.4288  # N137
$ws.Cells.Item(138, 8).Value = 4604.3267  # H138
$ws.Cells.Item(138, 9).Value = 912  # I138
$ws.Cells.Item(138, 10).Value = 5801.838  # J138
$ws.Cells.Item(138, 11).Value = 2736  # K138
$ws.Cells.Item(138, 12).Value = 17405.514  # L138
$ws.Cells.Item(138, 13).Value = 2404  # M138
$ws.Cells.Item(138, 14).Value = -27685.514  # N138

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(10, 8).Value = 13309.667  # H10
$ws.Cells.Item(10, 9).Value = 6004  # I10
$ws.Cells.Item(10, 10).Value = 16962.5  # J10
$ws.Cells.Item(10, 11).Value = 6004  # K10
$ws.Cells.Item(10, 12).Value = 16962.5  # L10
$ws.Cells.Item(10, 13).Value = -5834  # M10
$ws.Cells.Item(10, 14).Value = -17302.5  # N10
$ws.Cells.Item(32, 8).Value = 4219.4  # H32
$ws.Cells.Item(32, 9).Value = 3822.3962  # I32
$ws.Cells.Item(32, 10).Value = 7225.2856  # J32
$ws.Cells.Item(32, 11).Value = 3822.3962  # K32
$ws.Cells.Item(32, 12).Value = 7225.2856  # L32
$ws.Cells.Item(32, 13).Value = -3535.3962  # M32
$ws.Cells.Item(32, 14).Value = -7799.2856  # N32
$ws.Cells.Item(61, 8).Value = 1068.0189  # H61
$ws.Cells.Item(61, 9).Value = 827.63043  # I61
$ws.Cells.Item(61, 10).Value = 2647.7144  # J61
$ws.Cells.Item(61, 11).Value = 827.63043  # K61
$ws.Cells.Item(61, 12).Value = 2647.7144  # L61
$ws.Cells.Item(61, 13).Value = -615.63043  # M61
$ws.Cells.Item(61, 14).Value = -3071.7144  # N61
$ws.Cells.Item(74, 8).Value = 3421.7368  # H74
$ws.Cells.Item(74, 9).Value = 3381.2258  # I74
$ws.Cells.Item(74, 11).Value = 3381.2258  # K74
$ws.Cells.Item(74, 13).Value = -2507.2258  # M74
$ws.Cells.Item(77, 8).Value = 3421.7368  # H77
$ws.Cells.Item(77, 9).Value = 3381.2258  # I77
$ws.Cells.Item(77, 11).Value = 16906.129  # K77
$ws.Cells.Item(77, 13).Value = -12538.129  # M77
$ws.Cells.Item(102, 8).Value = 2177.7778  # H102
$ws.Cells.Item(102, 9).Value = 2100  # I102
$ws.Cells.Item(102, 10).Value = 2333.3333  # J102
$ws.Cells.Item(102, 11).Value = 2100  # K102
$ws.Cells.Item(102, 12).Value = 2333.3333  # L102
$ws.Cells.Item(102, 13).Value = -478  # M102
$ws.Cells.Item(102, 14).Value = -5577.3333  # N102
$ws.Cells.Item(103, 8).Value = 37400  # H103
$ws.Cells.Item(103, 10).Value = 37400  # J103
$ws.Cells.Item(103, 12).Value = 37400  # L103
$ws.Cells.Item(103, 14).Value = -39744  # N103
$ws.Cells.Item(110, 8).Value = 2405.7646  # H110
$ws.Cells.Item(110, 9).Value = 2352.6924  # I110
$ws.Cells.Item(110, 10).Value = 2578.25  # J110
$ws.Cells.Item(110, 11).Value = 2352.6924  # K110
$ws.Cells.Item(110, 12).Value = 2578.25  # L110
$ws.Cells.Item(110, 13).Value = -307.6923999999999  # M110
$ws.Cells.Item(110, 14).Value = -6668.25  # N110
$ws.Cells.Item(132, 8).Value = 2969.7942  # H132
$ws.Cells.Item(132, 9).Value = 1825.0435  # I132
$ws.Cells.Item(132, 10).Value = 5363.364  # J132
$ws.Cells.Item(132, 11).Value = 5475.1305  # K132
$ws.Cells.Item(132, 12).Value = 16090.092  # L132
$ws.Cells.Item(132, 13).Value = -2945.1305  # M132
$ws.Cells.Item(132, 14).Value = -21150.092  # N132
$ws.Cells.Item(136, 8).Value = 1068.0189  # H136
$ws.Cells.Item(136, 9).Value = 827.63043  # I136
$ws.Cells.Item(136, 10).Value = 2647.7144  # J136
$ws.Cells.Item(136, 11).Value = 2482.89129  # K136
$ws.Cells.Item(136, 12).Value = 7943.1432  # L136
$ws.Cells.Item(136, 13).Value = 67.10870999999997  # M136
$ws.Cells.Item(136, 14).Value = -13043.1432  # N136
$ws.Cells.Item(137, 8).Value = 30866.666  # H137
$ws.Cells.Item(137, 10).Value = 39800  # J137
$ws.Cells.Item(137, 12).Value = 39800  # L137
$ws.Cells.Item(137, 14).Value = -50000  # N137

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(59, 8).Value = 39450  # H59
$ws.Cells.Item(59, 10).Value = 39450  # J59
$ws.Cells.Item(59, 12).Value = 39450  # L59
$ws.Cells.Item(59, 14).Value = -41144  # N59
$ws.Cells.Item(137, 8).Value = 34180  # H137
$ws.Cells.Item(137, 10).Value = 37240  # J137
$ws.Cells.Item(137, 12).Value = 37240  # L137
$ws.Cells.Item(137, 14).Value = -47440  # N137

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(6, 8).Value = 4602010  # H6
$ws.Cells.Item(6, 10).Value = 10000  # J6
$ws.Cells.Item(6, 12).Value = 10000  # L6
$ws.Cells.Item(6, 14).Value = -10226  # N6
$ws.Cells.Item(31, 8).Value = 9806967  # H31
$ws.Cells.Item(31, 9).Value = 1754.9259  # I31
$ws.Cells.Item(31, 10).Value = 20837830  # J31
$ws.Cells.Item(31, 11).Value = 1754.9259  # K31
$ws.Cells.Item(31, 12).Value = 20837830  # L31
$ws.Cells.Item(31, 13).Value = -1459.9259  # M31
$ws.Cells.Item(31, 14).Value = -20838420  # N31
$ws.Cells.Item(34, 8).Value = 9806967  # H34
$ws.Cells.Item(34, 9).Value = 1754.9259  # I34
$ws.Cells.Item(34, 10).Value = 20837830  # J34
$ws.Cells.Item(34, 11).Value = 1754.9259  # K34
$ws.Cells.Item(34, 12).Value = 20837830  # L34
$ws.Cells.Item(34, 13).Value = -1552.9259  # M34
$ws.Cells.Item(34, 14).Value = -20838234  # N34
$ws.Cells.Item(58, 8).Value = 1563.8081  # H58
$ws.Cells.Item(58, 9).Value = 1459.5916  # I58
$ws.Cells.Item(58, 10).Value = 1828.0714  # J58
$ws.Cells.Item(58, 11).Value = 1459.5916  # K58
$ws.Cells.Item(58, 12).Value = 1828.0714  # L58
$ws.Cells.Item(58, 13).Value = -1256.5916  # M58
$ws.Cells.Item(58, 14).Value = -2234.0714  # N58
$ws.Cells.Item(134, 8).Value = 4170.8374  # H134
$ws.Cells.Item(134, 9).Value = 5197.75  # I134
$ws.Cells.Item(134, 10).Value = 2873.6843  # J134
$ws.Cells.Item(134, 11).Value = 15593.25  # K134
$ws.Cells.Item(134, 12).Value = 8621.052899999999  # L134
$ws.Cells.Item(134, 13).Value = -13058.25  # M134
$ws.Cells.Item(134, 14).Value = -13691.0529  # N134
$ws.Cells.Item(136, 8).Value = 1563.8081  # H136
$ws.Cells.Item(136, 9).Value = 1459.5916  # I136
$ws.Cells.Item(136, 10).Value = 1828.0714  # J136
$ws.Cells.Item(136, 11).Value = 4378.7748  # K136
$ws.Cells.Item(136, 12).Value = 5484.2142  # L136
$ws.Cells.Item(136, 13).Value = -1828.7748  # M136
$ws.Cells.Item(136, 14).Value = -10584.2142  # N136
$ws.Cells.Item(137, 8).Value = 49775  # H137
$ws.Cells.Item(137, 10).Value = 49775  # J137
$ws.Cells.Item(137, 12).Value = 49775  # L137
$ws.Cells.Item(137, 14).Value = -59975  # N137

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 1460.5172  # H5
$ws.Cells.Item(5, 9).Value = 357.57895  # I5
$ws.Cells.Item(5, 10).Value = 3556.1  # J5
$ws.Cells.Item(5, 11).Value = 1072.73685  # K5
$ws.Cells.Item(5, 12).Value = 10668.3  # L5
$ws.Cells.Item(5, 13).Value = -960.73685  # M5
$ws.Cells.Item(5, 14).Value = -10892.3  # N5
$ws.Cells.Item(122, 8).Value = 2432.6  # H122
$ws.Cells.Item(122, 9).Value = 471.73334  # I122
$ws.Cells.Item(122, 10).Value = 3413.0334  # J122
$ws.Cells.Item(122, 11).Value = 4245.60006  # K122
$ws.Cells.Item(122, 12).Value = 30717.3006  # L122
$ws.Cells.Item(122, 13).Value = -1795.60006  # M122
$ws.Cells.Item(122, 14).Value = -35617.3006  # N122
$ws.Cells.Item(131, 8).Value = 818.4820999999999  # H131
$ws.Cells.Item(131, 9).Value = 481.66666  # I131
$ws.Cells.Item(131, 10).Value = 882.9787  # J131
$ws.Cells.Item(131, 11).Value = 1444.99998  # K131
$ws.Cells.Item(131, 12).Value = 2648.9361  # L131
$ws.Cells.Item(131, 13).Value = 3595.00002  # M131
$ws.Cells.Item(131, 14).Value = -12728.9361  # N131
$ws.Cells.Item(135, 8).Value = 1460.5172  # H135
$ws.Cells.Item(135, 9).Value = 357.57895  # I135
$ws.Cells.Item(135, 10).Value = 3556.1  # J135
$ws.Cells.Item(135, 11).Value = 3218.21055  # K135
$ws.Cells.Item(135, 12).Value = 32004.9  # L135
$ws.Cells.Item(135, 13).Value = -683.2105500000002  # M135
$ws.Cells.Item(135, 14).Value = -37074.89999999999  # N135

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 2167.1035  # H102
$ws.Cells.Item(102, 9).Value = 1798  # I102
$ws.Cells.Item(102, 10).Value = 2621.3845  # J102
$ws.Cells.Item(102, 11).Value = 1798  # K102
$ws.Cells.Item(102, 12).Value = 2621.3845  # L102
$ws.Cells.Item(102, 13).Value = -176  # M102
$ws.Cells.Item(102, 14).Value = -5865.3845  # N102
$ws.Cells.Item(107, 8).Value = 4831531.5  # H107
$ws.Cells.Item(107, 9).Value = 303.33334  # I107
$ws.Cells.Item(107, 10).Value = 10101962  # J107
$ws.Cells.Item(107, 11).Value = 303.33334  # K107
$ws.Cells.Item(107, 12).Value = 10101962  # L107
$ws.Cells.Item(107, 13).Value = 1616.66666  # M107
$ws.Cells.Item(107, 14).Value = -10105802  # N107
$ws.Cells.Item(126, 8).Value = 5131.074  # H126
$ws.Cells.Item(126, 9).Value = 0  # I126
$ws.Cells.Item(126, 10).Value = 5131.074  # J126
$ws.Cells.Item(126, 11).Value = 0  # K126
$ws.Cells.Item(126, 12).Value = 15393.222  # L126
$ws.Cells.Item(126, 13).ClearContents()  # M126
$ws.Cells.Item(126, 14).Value = -20333.222  # N126
$ws.Cells.Item(137, 8).Value = 42790  # H137
$ws.Cells.Item(137, 10).Value = 42790  # J137
$ws.Cells.Item(137, 12).Value = 42790  # L137
$ws.Cells.Item(137, 14).Value = -52990  # N137

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(123, 8).Value = 39947.184  # H123
$ws.Cells.Item(123, 10).Value = 39947.184  # J123
$ws.Cells.Item(123, 12).Value = 39947.184  # L123
$ws.Cells.Item(123, 14).Value = -49747.184  # N123
$ws.Cells.Item(132, 8).Value = 8634.115  # H132
$ws.Cells.Item(132, 9).Value = 8717.151  # I132
$ws.Cells.Item(132, 11).Value = 26151.453  # K132
$ws.Cells.Item(132, 13).Value = -23621.453  # M132
$ws.Cells.Item(136, 8).Value = 2343.4468  # H136
$ws.Cells.Item(136, 9).Value = 1361.1351  # I136
$ws.Cells.Item(136, 10).Value = 5978  # J136
$ws.Cells.Item(136, 11).Value = 4083.4053  # K136
$ws.Cells.Item(136, 12).Value = 17934  # L136
$ws.Cells.Item(136, 13).Value = -1533.4053  # M136
$ws.Cells.Item(136, 14).Value = -23034  # N136

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(96, 8).Value = 8577856  # H96
$ws.Cells.Item(96, 9).Value = 5250001  # I96
$ws.Cells.Item(96, 10).Value = 11905712  # J96
$ws.Cells.Item(96, 11).Value = 5250001  # K96
$ws.Cells.Item(96, 12).Value = 11905712  # L96
$ws.Cells.Item(96, 13).Value = -5248628  # M96
$ws.Cells.Item(96, 14).Value = -11908458  # N96
$ws.Cells.Item(107, 8).Value = 854.1818  # H107
$ws.Cells.Item(107, 9).Value = 685.1429000000001  # I107
$ws.Cells.Item(107, 10).Value = 1150  # J107
$ws.Cells.Item(107, 11).Value = 2055.4287  # K107
$ws.Cells.Item(107, 12).Value = 3450  # L107
$ws.Cells.Item(107, 13).Value = -135.4287000000004  # M107
$ws.Cells.Item(107, 14).Value = -7290  # N107
$ws.Cells.Item(126, 8).Value = 267956.12  # H126
$ws.Cells.Item(126, 9).Value = 1210  # I126
$ws.Cells.Item(126, 10).Value = 628847.9399999999  # J126
$ws.Cells.Item(126, 11).Value = 3630  # K126
$ws.Cells.Item(126, 12).Value = 1886543.82  # L126
$ws.Cells.Item(126, 13).Value = -1160  # M126
$ws.Cells.Item(126, 14).Value = -1891483.82  # N126
$ws.Cells.Item(132, 8).Value = 7577474.5  # H132
$ws.Cells.Item(132, 9).Value = 690.9655  # I132
$ws.Cells.Item(132, 10).Value = 22225922  # J132
$ws.Cells.Item(132, 11).Value = 2072.8965  # K132
$ws.Cells.Item(132, 12).Value = 66677766  # L132
$ws.Cells.Item(132, 13).Value = 457.1035000000002  # M132
$ws.Cells.Item(132, 14).Value = -66682826  # N132
$ws.Cells.Item(136, 8).Value = 2003.0968  # H136
$ws.Cells.Item(136, 9).Value = 745.3946999999999  # I136
$ws.Cells.Item(136, 10).Value = 3994.4583  # J136
$ws.Cells.Item(136, 11).Value = 2236.1841  # K136
$ws.Cells.Item(136, 12).Value = 11983.3749  # L136
$ws.Cells.Item(136, 13).Value = 313.8159000000001  # M136
$ws.Cells.Item(136, 14).Value = -17083.3749  # N136
$ws.Cells.Item(139, 8).Value = 38642.69  # H139
$ws.Cells.Item(139, 10).Value = 38562.4  # J139
$ws.Cells.Item(139, 12).Value = 38562.4  # L139
$ws.Cells.Item(139, 14).Value = -48842.4  # N139
